# Update the torque summary table with the new data set.
# Values are kept as literal text (leading apostrophe forces text
# interpretation for the purely numeric-looking entries) so that the
# cell content/type matches the original sheet, where every data cell
# is stored as text rather than a numeric value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Applied Torque 230 -> 550
$ws.Range("A2").Value = "'550"
$ws.Range("B2").Value = "528.0 - 572.0"
$ws.Range("C2").Value = "'565.4"
$ws.Range("D2").Value = "'556.5"
$ws.Range("E2").Value = "'554.9"
$ws.Range("F2").Value = "'560.5"
$ws.Range("G2").Value = "'556.7"

# Row 3: Applied Torque 150 -> 350
$ws.Range("A3").Value = "'350"
$ws.Range("B3").Value = "336.0 - 364.0"
$ws.Range("C3").Value = "'355.7"
$ws.Range("D3").Value = "'342.1"
$ws.Range("E3").Value = "'343.7"
$ws.Range("F3").Value = "'341.1"
$ws.Range("G3").Value = "'339.9"

# Row 4: Applied Torque 70 -> 200
$ws.Range("A4").Value = "'200"
$ws.Range("B4").Value = "192.0 - 208.0"
$ws.Range("C4").Value = "'200.8"
$ws.Range("D4").Value = "'195.4"
$ws.Range("E4").Value = "'193.1"
$ws.Range("F4").Value = "'194.8"
$ws.Range("G4").Value = "'192.0"
